$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "56.619.57"
Set-TextCell "E2" "  -4.52%  "

Set-TextCell "D3" "2.359.94"
Set-TextCell "E3" "  -6.60%  "

Set-TextCell "E4" "  +0.04%  "

Set-TextCell "D5" "513.60"
Set-TextCell "E5" "  -4.27%  "

Set-TextCell "D6" "127.83"
Set-TextCell "E6" "  -6.01%  "

Set-TextCell "D7" "0.998"
Set-TextCell "E7" "  +0.32%  "

Set-TextCell "D8" "0.553"
Set-TextCell "E8" "  -2.47%  "

Set-TextCell "D9" "2.373.66"
Set-TextCell "E9" "  -6.01%  "

Set-TextCell "D10" "0.0956"
Set-TextCell "E10" "  -4.24%  "

Set-TextCell "E11" "  -2.00%  "

Set-TextCell "E12" "  -8.52%  "

Set-TextCell "D13" "0.316"
Set-TextCell "E13" "  -5.98%  "

Set-TextCell "D14" "2.779.00"
Set-TextCell "E14" "  -6.53%  "

Set-TextCell "D15" "56.511.68"
Set-TextCell "E15" "  -4.61%  "

Set-TextCell "D16" "21.43"
Set-TextCell "E16" "  -4.82%  "

Set-TextCell "D17" "0.0000131"
Set-TextCell "E17" "  -4.68%  "

Set-TextCell "D18" "2.373.55"
Set-TextCell "E18" "  -5.92%  "

Set-TextCell "D19" "10.28"
Set-TextCell "E19" "  -4.48%  "

Set-TextCell "D20" "309.56"
Set-TextCell "E20" "  -4.25%  "

Set-TextCell "D21" "4.01"
Set-TextCell "E21" "  -5.69%  "

Set-TextCell "D22" "6.07"
Set-TextCell "E22" "  -1.50%  "

Set-TextCell "D23" "0.998"
Set-TextCell "E23" "  +0.00%  "

Set-TextCell "D24" "64.59"
Set-TextCell "E24" "  -2.22%  "

Set-TextCell "D25" "0.997"
Set-TextCell "E25" "  -0.17%  "

Set-TextCell "D26" "0.389"
Set-TextCell "E26" "  -5.06%  "

Set-TextCell "D27" "2.461.96"
Set-TextCell "E27" "  -6.88%  "

Set-TextCell "D28" "0.153"
Set-TextCell "E28" "  -5.08%  "

Set-TextCell "D29" "7.16"
Set-TextCell "E29" "  -5.29%  "

Set-TextCell "D30" "172.62"
Set-TextCell "E30" "  -0.59%  "

Set-TextCell "E31" "  -5.35%  "

Set-TextCell "D32" "0.0₃0715"
Set-TextCell "E32" "  -7.05%  "

Set-TextCell "D33" "6.07"
Set-TextCell "E33" "  -5.15%  "

Set-TextCell "D34" "1.13"
Set-TextCell "E34" "  -7.76%  "

Set-TextCell "E35" "  -0.09%  "

Set-TextCell "D36" "0.995"
Set-TextCell "E36" "  -0.43%  "

Set-TextCell "D37" "17.56"
Set-TextCell "E37" "  -3.86%  "

Set-TextCell "E38" "  -6.59%  "

Set-TextCell "D39" "3.72"
Set-TextCell "E39" "  -7.70%  "

Set-TextCell "B40" "SuiNetwork"
Set-TextCell "C40" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell "D40" "0.799"
Set-TextCell "E40" "  +0.76%  "

Set-TextCell "B41" "OKB"
Set-TextCell "C41" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D41" "35.46"
Set-TextCell "E41" "  -3.71%  "

Set-TextCell "D42" "1.42"
Set-TextCell "E42" "  -7.02%  "

Set-TextCell "D43" "3.32"
Set-TextCell "E43" "  -5.61%  "

Set-TextCell "D44" "4.88"
Set-TextCell "E44" "  -4.92%  "

Set-TextCell "D45" "123.04"
Set-TextCell "E45" "  -7.12%  "

Set-TextCell "D46" "0.568"
Set-TextCell "E46" "  -5.53%  "

Set-TextCell "D47" "252.37"
Set-TextCell "E47" "  -10.49%  "

Set-TextCell "D48" "0.0906"
Set-TextCell "E48" "  -2.67%  "

Set-TextCell "D49" "0.0487"
Set-TextCell "E49" "  -5.05%  "

Set-TextCell "D50" "0.0207"
Set-TextCell "E50" "  -6.33%  "

Set-TextCell "D51" "16.66"
Set-TextCell "E51" "  -6.61%  "
